$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value() = '41.012.82'
$ws.Range('E2').Value() = '  -1.13%  '
$ws.Range('D3').Value() = '2.174.08'
$ws.Range('E3').Value() = '  -2.13%  '
$ws.Range('D4').Value() = '0.999'
$ws.Range('E4').Value() = '  -0.18%  '
$ws.Range('D5').Value() = '247.08'
$ws.Range('E5').Value() = '  -1.15%  '
$ws.Range('D6').Value() = '0.615'
$ws.Range('E6').Value() = '  -1.58%  '
$ws.Range('D7').Value() = '66.49'
$ws.Range('E7').Value() = '  -5.93%  '
$ws.Range('E8').Value() = '  -0.12%  '
$ws.Range('D9').Value() = '0.568'
$ws.Range('E9').Value() = '  +0.05%  '
$ws.Range('D10').Value() = '60.49'
$ws.Range('E10').Value() = '  +3.07%  '
$ws.Range('D11').Value() = '0.0929'
$ws.Range('E11').Value() = '  -2.85%  '
$ws.Range('D12').Value() = '35.78'
$ws.Range('E12').Value() = '  -14.14%  '
$ws.Range('D13').Value() = '0.104'
$ws.Range('E13').Value() = '  -1.65%  '
$ws.Range('D14').Value() = '6.81'
$ws.Range('E14').Value() = '  -2.46%  '
$ws.Range('D15').Value() = '2.477.60'
$ws.Range('E15').Value() = '  -3.04%  '
$ws.Range('D16').Value() = '0.856'
$ws.Range('E16').Value() = '  +0.41%  '
$ws.Range('D17').Value() = '14.30'
$ws.Range('E17').Value() = '  -3.67%  '
$ws.Range('D18').Value() = '2.181.21'
$ws.Range('E18').Value() = '  -1.76%  '
$ws.Range('D19').Value() = '40.944.74'
$ws.Range('E19').Value() = '  -1.26%  '
$ws.Range('D20').Value() = '0.0₃0940'
$ws.Range('E20').Value() = '  -2.20%  '
$ws.Range('D21').Value() = '6.09'
$ws.Range('E21').Value() = '  -1.65%  '
$ws.Range('D22').Value() = '71.44'
$ws.Range('E22').Value() = '  -1.55%  '
$ws.Range('D23').Value() = '229.99'
$ws.Range('E23').Value() = '  -1.66%  '
$ws.Range('D24').Value() = '2.07'
$ws.Range('E24').Value() = '  -7.53%  '
$ws.Range('E25').Value() = '  +0.15%  '
$ws.Range('D26').Value() = '11.36'
$ws.Range('E26').Value() = '  +7.29%  '
$ws.Range('D27').Value() = '3.68'
$ws.Range('E27').Value() = '  -3.90%  '
$ws.Range('D28').Value() = '2.43'
$ws.Range('E28').Value() = '  -2.11%  '
$ws.Range('E29').Value() = '  -5.56%  '
$ws.Range('D30').Value() = '168.72'
$ws.Range('E30').Value() = '  -1.31%  '
$ws.Range('D31').Value() = '2.01'
$ws.Range('E31').Value() = '  -8.16%  '
$ws.Range('D32').Value() = '20.21'
$ws.Range('E32').Value() = '  -1.74%  '
$ws.Range('D33').Value() = '0.121'
$ws.Range('E33').Value() = '  +0.96%  '
$ws.Range('D34').Value() = '5.65'
$ws.Range('E34').Value() = '  +2.08%  '
$ws.Range('E35').Value() = '  +4.83%  '
$ws.Range('E36').Value() = '  -2.26%  '
$ws.Range('D37').Value() = '4.56'
$ws.Range('E37').Value() = '  -2.44%  '
$ws.Range('D38').Value() = '4.04'
$ws.Range('E38').Value() = '  +2.38%  '
$ws.Range('D39').Value() = '24.28'
$ws.Range('E39').Value() = '  -6.82%  '
$ws.Range('D40').Value() = '0.0303'
$ws.Range('E40').Value() = '  +4.79%  '
$ws.Range('E41').Value() = '  -4.32%  '
$ws.Range('D42').Value() = '5.47'
$ws.Range('E42').Value() = '  -8.01%  '
$ws.Range('D43').Value() = '4.86'
$ws.Range('E43').Value() = '  -0.30%  '
$ws.Range('D44').Value() = '11.32'
$ws.Range('E44').Value() = '  -5.41%  '
$ws.Range('D45').Value() = '60.59'
$ws.Range('E45').Value() = '  -11.99%  '
$ws.Range('D46').Value() = '0.193'
$ws.Range('E46').Value() = '  -6.74%  '
$ws.Range('B47').Value() = 'FraxShare'
$ws.Range('C47').Value() = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value() = '8.45'
$ws.Range('E47').Value() = '  -3.09%  '
$ws.Range('B48').Value() = 'Cronos'
$ws.Range('C48').Value() = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value() = '0.0991'
$ws.Range('E48').Value() = '  -2.22%  '
$ws.Range('B49').Value() = 'BinanceUSD'
$ws.Range('C49').Value() = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D49').Value() = '1.00'
$ws.Range('E49').Value() = '  +0.00%  '
$ws.Range('E50').Value() = '  -0.32%  '
$ws.Range('E51').Value() = '  -3.09%  '
